$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.935.08'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '2.353.90'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.21'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.669'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.34%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '74.53'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.69%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.597'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.52%  '
$ws.Range('E10').Value = '  +0.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '59.80'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.43%  '
$ws.Range('E12').Value = '  +0.41%  '
$ws.Range('E13').Value = '  +1.13%  '
$ws.Range('E14').Value = '  -1.37%  '
$ws.Range('D15').Value = '2.699.63'
$ws.Range('E15').Value = '  -0.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '16.23'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.907'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('D18').Value = '2.350.01'
$ws.Range('E18').Value = '  +1.30%  '
$ws.Range('D19').Value = '43.886.69'
$ws.Range('E19').Value = '  +0.02%  '
$ws.Range('E20').Value = '  +1.53%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.67'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.75%  '
$ws.Range('B22').Value = 'Litecoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '78.32'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '253.43'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.91%  '
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('E25').Value = '  +1.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.85'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.51'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.49'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.76%  '
$ws.Range('E29').Value = '  +1.48%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '176.49'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.32'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.44%  '
$ws.Range('E32').Value = '  +0.83%  '
$ws.Range('E33').Value = '  -1.74%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0751'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.09'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.40'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.84'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.95%  '
$ws.Range('B38').Value = 'THORChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.43'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.91%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.39'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.58'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +17.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0272'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '65.15'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +16.73%  '
$ws.Range('E43').Value = '  +1.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.01'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.89%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.201'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.43%  '
$ws.Range('E46').Value = '  -3.18%  '
$ws.Range('E47').Value = '  +0.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.23'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.43'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.16'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.39%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '98.68'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.10%  '
